$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before CO, shifting existing CO column (and its data) to CP
$ws.Columns("CO:CO").Insert()

# Fill header and "Faixa Etaria" (age range) values for the new CO column
$ws.Range("CO1").Value = "Faixa Etária"
$ws.Range("CO2").Value = "entre 20 e 25"
$ws.Range("CO3").Value = "Entre 15 e 20"
$ws.Range("CO4").Value = "Entre 15 e 20"
$ws.Range("CO5").Value = "Entre 15 e 20"
$ws.Range("CO6").Value = "entre 30 e 35"
$ws.Range("CO7").Value = "entre 20 e 25"
$ws.Range("CO8").Value = "Entre 15 e 20"
$ws.Range("CO9").Value = "Entre 15 e 20"
$ws.Range("CO10").Value = "Entre 15 e 20"
$ws.Range("CO11").Value = "entre 25 e 30"
$ws.Range("CO12").Value = "entre 20 e 25"
$ws.Range("CO13").Value = "entre 30 e 35"
$ws.Range("CO14").Value = "entre 35 e 40"
$ws.Range("CO15").Value = "entre 35 e 40"
$ws.Range("CO16").Value = "Acima de 40"
$ws.Range("CO17").Value = "entre 30 e 35"
$ws.Range("CO18").Value = "entre 20 e 25"
$ws.Range("CO19").Value = "entre 25 e 30"
$ws.Range("CO20").Value = "entre 20 e 25"
$ws.Range("CO21").Value = "Entre 15 e 20"
$ws.Range("CO22").Value = "Entre 15 e 20"
$ws.Range("CO23").Value = "Entre 15 e 20"
$ws.Range("CO24").Value = "entre 20 e 25"
$ws.Range("CO25").Value = "entre 25 e 30"
$ws.Range("CO26").Value = "Entre 15 e 20"
$ws.Range("CO27").Value = "Entre 15 e 20"
$ws.Range("CO28").Value = "Entre 15 e 20"
$ws.Range("CO29").Value = "entre 25 e 30"
$ws.Range("CO30").Value = "Entre 15 e 20"
$ws.Range("CO31").Value = "entre 20 e 25"
$ws.Range("CO32").Value = "entre 30 e 35"
$ws.Range("CO33").Value = "Entre 15 e 20"
$ws.Range("CO34").Value = "Entre 15 e 20"
$ws.Range("CO35").Value = "entre 30 e 35"
$ws.Range("CO36").Value = "entre 20 e 25"
$ws.Range("CO37").Value = "entre 20 e 25"
$ws.Range("CO38").Value = "entre 30 e 35"
$ws.Range("CO39").Value = "entre 20 e 25"
$ws.Range("CO40").Value = "entre 20 e 25"
$ws.Range("CO41").Value = "entre 35 e 40"
$ws.Range("CO42").Value = "Entre 15 e 20"
$ws.Range("CO43").Value = "entre 30 e 35"
$ws.Range("CO44").Value = "entre 30 e 35"
$ws.Range("CO45").Value = "entre 35 e 40"
$ws.Range("CO46").Value = "Entre 15 e 20"
$ws.Range("CO47").Value = "Entre 15 e 20"
$ws.Range("CO48").Value = "entre 35 e 40"
$ws.Range("CO49").Value = "entre 30 e 35"
$ws.Range("CO50").Value = "entre 20 e 25"
$ws.Range("CO51").Value = "entre 30 e 35"
$ws.Range("CO52").Value = "entre 20 e 25"
$ws.Range("CO53").Value = "entre 30 e 35"
$ws.Range("CO54").Value = "entre 20 e 25"
$ws.Range("CO55").Value = "entre 20 e 25"
$ws.Range("CO56").Value = "Acima de 40"
$ws.Range("CO57").Value = "Entre 15 e 20"
$ws.Range("CO58").Value = "entre 20 e 25"
$ws.Range("CO59").Value = "entre 30 e 35"
$ws.Range("CO60").Value = "Entre 15 e 20"
$ws.Range("CO61").Value = "Entre 15 e 20"
$ws.Range("CO62").Value = "Acima de 40"
$ws.Range("CO63").Value = "entre 25 e 30"
$ws.Range("CO64").Value = "entre 25 e 30"
$ws.Range("CO65").Value = "entre 25 e 30"
$ws.Range("CO66").Value = "Entre 15 e 20"
$ws.Range("CO67").Value = "entre 30 e 35"

Write-Host "Done filling CO column"
